{"js": "// Apply exact text replacements for the three-digit / one-digit division answers table.\n// Each (old, new) pair is unique within the document, so a plain case-sensitive\n// search-and-replace on context.document.body is sufficient and keeps the original\n// run formatting (rFonts/sz) intact since insertText(..., replace) rewrites the run's\n// text in place rather than creating a new run.\nconst replacements = [\n  [\"321\u00f76=53, 3\", \"619\u00f73=206, 1\"],\n  [\"212\u00f76=35, 2\", \"943\u00f76=157, 1\"],\n  [\"352\u00f72=176, 0\", \"925\u00f78=115, 5\"],\n  [\"358\u00f76=59, 4\", \"976\u00f72=488, 0\"],\n  [\"571\u00f75=114, 1\", \"612\u00f73=204, 0\"],\n  [\"169\u00f76=28, 1\", \"881\u00f79=97, 8\"],\n  [\"474\u00f72=237, 0\", \"684\u00f72=342, 0\"],\n  [\"250\u00f77=35, 5\", \"376\u00f75=75, 1\"],\n  [\"848\u00f72=424, 0\", \"163\u00f79=18, 1\"],\n  [\"962\u00f76=160, 2\", \"334\u00f72=167, 0\"],\n  [\"753\u00f74=188, 1\", \"713\u00f75=142, 3\"],\n  [\"812\u00f74=203, 0\", \"365\u00f78=45, 5\"],\n  [\"191\u00f75=38, 1\", \"408\u00f73=136, 0\"],\n  [\"313\u00f77=44, 5\", \"971\u00f73=323, 2\"],\n  [\"109\u00f73=36, 1\", \"171\u00f74=42, 3\"],\n  [\"648\u00f73=216, 0\", \"725\u00f77=103, 4\"],\n  [\"873\u00f77=124, 5\", \"664\u00f77=94, 6\"],\n  [\"124\u00f73=41, 1\", \"655\u00f79=72, 7\"],\n  [\"374\u00f78=46, 6\", \"858\u00f73=286, 0\"],\n  [\"725\u00f75=145, 0\", \"400\u00f72=200, 0\"],\n  [\"218\u00f76=36, 2\", \"658\u00f78=82, 2\"],\n  [\"573\u00f77=81, 6\", \"153\u00f72=76, 1\"],\n  [\"392\u00f75=78, 2\", \"631\u00f76=105, 1\"],\n  [\"213\u00f73=71, 0\", \"268\u00f76=44, 4\"],\n  [\"830\u00f77=118, 4\", \"940\u00f78=117, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\" but found ${results.items.length}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply exact text replacements for the three-digit / one-digit division answers table.\n# Each (old, new) pair is unique within the document, so Find/Replace against the whole\n# document Range is sufficient; Word's Find.Execute rewrites the w:t text of the matched\n# run in place, leaving run formatting (rFonts/sz) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('321\u00f76=53, 3', '619\u00f73=206, 1'),\n    @('212\u00f76=35, 2', '943\u00f76=157, 1'),\n    @('352\u00f72=176, 0', '925\u00f78=115, 5'),\n    @('358\u00f76=59, 4', '976\u00f72=488, 0'),\n    @('571\u00f75=114, 1', '612\u00f73=204, 0'),\n    @('169\u00f76=28, 1', '881\u00f79=97, 8'),\n    @('474\u00f72=237, 0', '684\u00f72=342, 0'),\n    @('250\u00f77=35, 5', '376\u00f75=75, 1'),\n    @('848\u00f72=424, 0', '163\u00f79=18, 1'),\n    @('962\u00f76=160, 2', '334\u00f72=167, 0'),\n    @('753\u00f74=188, 1', '713\u00f75=142, 3'),\n    @('812\u00f74=203, 0', '365\u00f78=45, 5'),\n    @('191\u00f75=38, 1', '408\u00f73=136, 0'),\n    @('313\u00f77=44, 5', '971\u00f73=323, 2'),\n    @('109\u00f73=36, 1', '171\u00f74=42, 3'),\n    @('648\u00f73=216, 0', '725\u00f77=103, 4'),\n    @('873\u00f77=124, 5', '664\u00f77=94, 6'),\n    @('124\u00f73=41, 1', '655\u00f79=72, 7'),\n    @('374\u00f78=46, 6', '858\u00f73=286, 0'),\n    @('725\u00f75=145, 0', '400\u00f72=200, 0'),\n    @('218\u00f76=36, 2', '658\u00f78=82, 2'),\n    @('573\u00f77=81, 6', '153\u00f72=76, 1'),\n    @('392\u00f75=78, 2', '631\u00f76=105, 1'),\n    @('213\u00f73=71, 0', '268\u00f76=44, 4'),\n    @('830\u00f77=118, 4', '940\u00f78=117, 4'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord (text includes \u00f7 / = / , so \"word\" boundaries don't apply)\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
